# Commit: "Use different months for summer for bsth, os, cs"
#
# The Climate sheet gains a new "app.mthd.wthr" column (C) so that weather
# timing data can be specified per application method, plus a "notes"
# column (G). The old single "Sommer/Summer" row is split into three rows
# (Trailing hose / Open slot injection / Closed slot injection), each with
# its own climate figures and an explanatory note. The other month rows get
# "All" in the new app.mthd.wthr column.

$wb = $excel.ActiveWorkbook
$climate = $wb.Worksheets.Item("Climate")
$application = $wb.Worksheets.Item("Application")

# --- Climate sheet -------------------------------------------------------

# Insert a new column C ("app.mthd.wthr"); old C/D/E (air.temp/wind.2m/
# rain.rate) shift right to D/E/F.
$climate.Range("C1").EntireColumn.Insert()

# Row 1 - headers
$climate.Range("C1").Value = "app.mthd.wthr"
$climate.Range("G1").Value = "notes"

# Row 2 - Marts / March
$climate.Range("C2").Value = "All"

# Row 3 - April / April
$climate.Range("C3").Value = "All"

# Row 4 - Maj / May
$climate.Range("C4").Value = "All"

# Row 5 - Sommer / Summer -> Trailing hose
$climate.Range("C5").Value = "Trailing hose"
$climate.Range("G5").Value = "For trailing hose, 6-8."

# Row 6 - Efterår / Autumn
$climate.Range("C6").Value = "All"
$climate.Range("G6").Value = "9 (September)"

# Pre-seed rows 7/8 formatting (style only, no values) by copying row 6's
# cell style down -- gives the new rows the same style index ("s=1") as
# every other data cell in the sheet.
$climate.Range("A6:G6").Copy()
$climate.Range("A7:G7").PasteSpecial(-4122)
$climate.Range("A6:G6").Copy()
$climate.Range("A8:G8").PasteSpecial(-4122)

# Row 7 (new) - Sommer / Summer -> Open slot injection
$climate.Range("A7").Value = "Sommer"
$climate.Range("B7").Value = "Summer"
$climate.Range("C7").Value = "Open slot injection"
$climate.Range("D7").Value = 15.75
$climate.Range("E7").Value = 3.2775
$climate.Range("F7").Value = 0.09
$climate.Range("G7").Value = "Summer-grass, for open slot injection, 5-8."

# Row 8 (new) - Sommer / Summer -> Closed slot injection
$climate.Range("A8").Value = "Sommer"
$climate.Range("B8").Value = "Summer"
$climate.Range("C8").Value = "Closed slot injection"
$climate.Range("D8").Value = 17.55
$climate.Range("E8").Value = 3.105
$climate.Range("F8").Value = 0.09
$climate.Range("G8").Value = "Summer, before winter rapeseed, for closed slot injection, 7-8."

# Column widths (B widens for the longer method text, new C and G columns
# get their own widths to fit the method / notes text).
$climate.Columns.Item(2).ColumnWidth = 27.571428571428573
$climate.Columns.Item(3).ColumnWidth = 16.571428571428573
$climate.Columns.Item(7).ColumnWidth = 51.285714285714285

# --- Application sheet selection/view ------------------------------------
# The previously-active tab (Application) loses the selection, moving its
# remembered cell cursor down past the data rows.
$application.Activate()
$application.Range("B12").Select()

# --- Climate becomes the active tab/selection -----------------------------
$climate.Activate()
$climate.Range("B9").Select()
